# "Generate Report for Handback"
#
# 1e7f6c2f-...md and a8817da4-...md move from "Ready for handoff" to
# "Handed back: in sync with en-US" on the Overview sheet, and on each of
# the zh-cn / de-de detail sheets they pick up a "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime". Rows are reordered
# so the freshly-handed-back files sort to the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Canonical source-of-truth hyperlink targets (unchanged by this edit -
# only which cell they are attached to, and their display text, moves).
# ---------------------------------------------------------------------
$urls = @{
    "a4934d57" = @{
        "md"       = "https://github.com/OpenLocalizationTest/oltest/blob/83bfee64a72118cb35b9bfdf7f456f93197dec83/e2e/a4934d57-5e8b-43ec-91f3-96ee9fc97e87.md"
        "zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9e9d227299e76b6660290b26d4ef9bfbd6d099b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a4934d57-5e8b-43ec-91f3-96ee9fc97e87.c20fd6aac3930c988d65488fb90c05fd0ad0e795.zh-cn.xlf"
        "de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7a5d20749a0f94543f856d49e3a6da67199a147f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a4934d57-5e8b-43ec-91f3-96ee9fc97e87.c20fd6aac3930c988d65488fb90c05fd0ad0e795.de-de.xlf"
    }
    "1e7f6c2f" = @{
        "md"       = "https://github.com/OpenLocalizationTest/oltest/blob/1ba53d3b0ebf937094b5e2fa8c63b96f2cc0cdb3/e2e/1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"
        "zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4ce9c5776d3bb1443258d88803b7b508d6abbee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.zh-cn.xlf"
        "de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/477ed96c4f62772c29258a5ed8233b81f56ab2d8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.de-de.xlf"
    }
    "a8817da4" = @{
        "md"       = "https://github.com/OpenLocalizationTest/oltest/blob/1ba53d3b0ebf937094b5e2fa8c63b96f2cc0cdb3/e2e/a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"
        "zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4ce9c5776d3bb1443258d88803b7b508d6abbee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.zh-cn.xlf"
        "de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/477ed96c4f62772c29258a5ed8233b81f56ab2d8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.de-de.xlf"
    }
    "d5ef8c1c" = @{
        "md"       = "https://github.com/OpenLocalizationTest/oltest/blob/cc4f896aaf60fe7253d8128a85fa11b1def1f3a3/e2e/d5ef8c1c-1440-43a9-9dca-75493d500fed.md"
        "zh-cn.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a4ce9c5776d3bb1443258d88803b7b508d6abbee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d5ef8c1c-1440-43a9-9dca-75493d500fed.d3a9915bbd5e008ccbcdec6266db7d24b217cd18.zh-cn.xlf"
        "de-de.xlf" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/477ed96c4f62772c29258a5ed8233b81f56ab2d8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d5ef8c1c-1440-43a9-9dca-75493d500fed.d3a9915bbd5e008ccbcdec6266db7d24b217cd18.de-de.xlf"
    }
}

function Clear-AllHyperlinks($ws) {
    while ($ws.Hyperlinks.Count -gt 0) {
        foreach ($h in $ws.Hyperlinks) {
            $h.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: reorder rows (handed-back files first) and flip their
# status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Clear-AllHyperlinks $wsOverview

$overviewRows = @(
    @{ id = "1e7f6c2f"; name = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"; status = "Handed back: in sync with en-US"; date = "2016-45-18 03:45:49" },
    @{ id = "a8817da4"; name = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"; status = "Handed back: in sync with en-US"; date = "2016-45-18 03:45:49" },
    @{ id = "a4934d57"; name = "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.md"; status = "In Translation";                 date = "2016-44-18 03:44:44" },
    @{ id = "d5ef8c1c"; name = "d5ef8c1c-1440-43a9-9dca-75493d500fed.md"; status = "Ready for handoff";              date = "2016-45-18 03:45:49" }
)

$row = 2
foreach ($r in $overviewRows) {
    $wsOverview.Cells.Item($row, 1).Value = $r.name
    $wsOverview.Cells.Item($row, 2).Value = $r.status
    $wsOverview.Cells.Item($row, 3).Value = $r.status
    $wsOverview.Cells.Item($row, 4).Value = $r.date
    $row = $row + 1
}

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urls["1e7f6c2f"]["md"], "", "", "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urls["a8817da4"]["md"], "", "", "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $urls["a4934d57"]["md"], "", "", "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), $urls["d5ef8c1c"]["md"], "", "", "d5ef8c1c-1440-43a9-9dca-75493d500fed.md")

# ---------------------------------------------------------------------
# Per-locale detail sheets (zh-cn / de-de): reorder rows, flip status,
# and populate Latest Target File / Latest Handback File / Latest
# Handback DateTime (columns F/G/H) for the two handed-back files.
# ---------------------------------------------------------------------
function Update-LocaleSheet($ws, $locale, $rows) {
    Clear-AllHyperlinks $ws

    $row = 2
    foreach ($r in $rows) {
        $ws.Cells.Item($row, 1).Value = $r.name            # A Source File Name
        $ws.Cells.Item($row, 2).Value = ".md"               # B File Extension
        $ws.Cells.Item($row, 3).Value = $r.status            # C Status
        $ws.Cells.Item($row, 4).Value = $r.handoffFile        # D Latest Handoff File
        $ws.Cells.Item($row, 5).Value = $r.handoffDate         # E Latest Handoff Datetime
        if ($r.targetFile) {
            $ws.Cells.Item($row, 6).Value = $r.targetFile       # F Latest Target File
            $ws.Cells.Item($row, 7).Value = $r.handbackFile      # G Latest Handback File
        }
        $ws.Cells.Item($row, 8).Value = $r.handbackDate          # H Latest Handback DateTime
        $ws.Cells.Item($row, 9).Value = "Include"                 # I Handoff Reason
        $row = $row + 1
    }

    $row = 2
    foreach ($r in $rows) {
        $aAddr = "A" + $row
        $bAddr = "B" + $row
        $dAddr = "D" + $row
        $ws.Hyperlinks.Add($ws.Range($aAddr), $urls[$r.id]["md"], "", "", $r.name)
        $ws.Hyperlinks.Add($ws.Range($bAddr), $urls[$r.id]["md"], "", "", ".md")
        $ws.Hyperlinks.Add($ws.Range($dAddr), $urls[$r.id][$locale + ".xlf"], "", "", $r.handoffFile)
        if ($r.targetFile) {
            $fAddr = "F" + $row
            $gAddr = "G" + $row
            $ws.Hyperlinks.Add($ws.Range($fAddr), $urls[$r.id]["md"], "", "", $r.targetFile)
            $ws.Hyperlinks.Add($ws.Range($gAddr), $urls[$r.id][$locale + ".xlf"], "", "", $r.handbackFile)
        }
        $row = $row + 1
    }
}

$zhRows = @(
    @{ id = "1e7f6c2f"; name = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"; status = "Handed back: in sync with en-US";
       handoffFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.zh-cn.xlf"; handoffDate = "2016-03-18 03:45:46";
       targetFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"; handbackFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.zh-cn.xlf"; handbackDate = "2016-03-18 03:46:09" },
    @{ id = "a8817da4"; name = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"; status = "Handed back: in sync with en-US";
       handoffFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.zh-cn.xlf"; handoffDate = "2016-03-18 03:45:46";
       targetFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"; handbackFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.zh-cn.xlf"; handbackDate = "2016-03-18 03:46:09" },
    @{ id = "a4934d57"; name = "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.md"; status = "In Translation";
       handoffFile = "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.c20fd6aac3930c988d65488fb90c05fd0ad0e795.zh-cn.xlf"; handoffDate = "2016-03-18 03:44:41";
       targetFile = $null; handbackFile = $null; handbackDate = "0001-01-01 00:00:00" },
    @{ id = "d5ef8c1c"; name = "d5ef8c1c-1440-43a9-9dca-75493d500fed.md"; status = "Ready for handoff";
       handoffFile = "d5ef8c1c-1440-43a9-9dca-75493d500fed.d3a9915bbd5e008ccbcdec6266db7d24b217cd18.zh-cn.xlf"; handoffDate = "2016-03-18 03:45:46";
       targetFile = $null; handbackFile = $null; handbackDate = "0001-01-01 00:00:00" }
)

$deRows = @(
    @{ id = "1e7f6c2f"; name = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"; status = "Handed back: in sync with en-US";
       handoffFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.de-de.xlf"; handoffDate = "2016-03-18 03:45:49";
       targetFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.md"; handbackFile = "1e7f6c2f-2faa-435f-939e-4ff5a4405080.5b62481d29727a02c4e9a9b0ca71e8731c76a1de.de-de.xlf"; handbackDate = "2016-03-18 03:46:14" },
    @{ id = "a8817da4"; name = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"; status = "Handed back: in sync with en-US";
       handoffFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.de-de.xlf"; handoffDate = "2016-03-18 03:45:49";
       targetFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.md"; handbackFile = "a8817da4-c4de-4d9a-8d12-d1e1d5844ebe.abb1196c605baac541d6d8558eae011236c19d95.de-de.xlf"; handbackDate = "2016-03-18 03:46:14" },
    @{ id = "a4934d57"; name = "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.md"; status = "In Translation";
       handoffFile = "a4934d57-5e8b-43ec-91f3-96ee9fc97e87.c20fd6aac3930c988d65488fb90c05fd0ad0e795.de-de.xlf"; handoffDate = "2016-03-18 03:44:44";
       targetFile = $null; handbackFile = $null; handbackDate = "0001-01-01 00:00:00" },
    @{ id = "d5ef8c1c"; name = "d5ef8c1c-1440-43a9-9dca-75493d500fed.md"; status = "Ready for handoff";
       handoffFile = "d5ef8c1c-1440-43a9-9dca-75493d500fed.d3a9915bbd5e008ccbcdec6266db7d24b217cd18.de-de.xlf"; handoffDate = "2016-03-18 03:45:49";
       targetFile = $null; handbackFile = $null; handbackDate = "0001-01-01 00:00:00" }
)

Update-LocaleSheet $wb.Worksheets.Item("zh-cn") "zh-cn" $zhRows
Update-LocaleSheet $wb.Worksheets.Item("de-de") "de-de" $deRows

Write-Host "Handback report generated."
